$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.783.71"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").Value = "2.673.29"
$ws.Range("E3").Value = "  -0.62%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.51"
$ws.Range("E5").Value = "  -1.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "157.63"
$ws.Range("E6").Value = "  -0.58%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.613"
$ws.Range("E8").Value = "  +3.90%  "
$ws.Range("E9").Value = "  +2.70%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.400"
$ws.Range("E10").Value = "  -0.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.83"
$ws.Range("E11").Value = "  -3.35%  "
$ws.Range("E12").Value = "  -0.23%  "
$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "29.08"
$ws.Range("E13").Value = "  -3.59%  "
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000200"
$ws.Range("E14").Value = "  -4.52%  "
$ws.Range("D15").Value = "3.155.43"
$ws.Range("E15").Value = "  -0.53%  "
$ws.Range("D16").Value = "65.669.89"
$ws.Range("E16").Value = "  -0.35%  "
$ws.Range("D17").Value = "2.670.53"
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.73"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.80"
$ws.Range("E19").Value = "  -1.82%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.51"
$ws.Range("E20").Value = "  -3.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "352.68"
$ws.Range("E21").Value = "  -1.75%  "
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.53"
$ws.Range("E23").Value = "  -2.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000115"
$ws.Range("E24").Value = "  +1.89%  "
$ws.Range("E25").Value = "  -1.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.68"
$ws.Range("E26").Value = "  +2.48%  "
$ws.Range("E27").Value = "  -3.46%  "
$ws.Range("E28").Value = "  -3.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.03"
$ws.Range("E29").Value = "  -3.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  -0.14%  "
$ws.Range("E31").Value = "  -3.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "531.85"
$ws.Range("E32").Value = "  -1.34%  "
$ws.Range("E33").Value = "  -0.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.49"
$ws.Range("E34").Value = "  -3.17%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.47"
$ws.Range("E35").Value = "  -0.48%  "
$ws.Range("E36").Value = "  -2.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.59"
$ws.Range("E37").Value = "  -1.02%  "
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "157.72"
$ws.Range("E39").Value = "  -4.19%  "
$ws.Range("E40").Value = "  -2.73%  "
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "163.35"
$ws.Range("E42").Value = "  -2.82%  "
$ws.Range("E43").Value = "  -1.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.37"
$ws.Range("E44").Value = "  +2.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0614"
$ws.Range("E45").Value = "  -3.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.82"
$ws.Range("E46").Value = "  -4.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0258"
$ws.Range("E47").Value = "  -3.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.640"
$ws.Range("E48").Value = "  -2.59%  "
$ws.Range("E49").Value = "  +8.92%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "20.03"
$ws.Range("E50").Value = "  -3.92%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0990"
$ws.Range("E51").Value = "  -0.49%  "
